# Update cryptos list data: apply per-cell text changes from the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-looking numeric strings (e.g. "6.50", "317.84")
# that Excel would otherwise auto-convert to Number on assignment, and multi-dot
# "thousands" strings (e.g. "41.622.74") that must stay literal text too. Force the
# whole Price column to Text format before writing, then drop back to the default
# "Normal" style afterwards so cell styling matches the original (unstyled) cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '41.622.74'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '2.473.51'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '317.84'
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('D6').Value = '92.17'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  +1.52%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.0861'
$ws.Range('E10').Value = '  +10.06%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').Value = '33.07'
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D13').Value = '2.854.50'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('E15').Value = '  -2.99%  '
$ws.Range('D16').Value = '2.467.66'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '0.789'
$ws.Range('E17').Value = '  +3.17%  '
$ws.Range('D18').Value = '41.570.85'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').Value = '6.50'
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').Value = '71.29'
$ws.Range('E21').Value = '  -0.64%  '
$ws.Range('D22').Value = '11.34'
$ws.Range('E22').Value = '  +1.10%  '
$ws.Range('D23').Value = '239.69'
$ws.Range('E23').Value = '  +1.41%  '
$ws.Range('D24').Value = '2.74'
$ws.Range('E24').Value = '  +1.36%  '
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '24.73'
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E28').Value = '  +3.90%  '
$ws.Range('D29').Value = '9.86'
$ws.Range('E29').Value = '  +2.53%  '
$ws.Range('D30').Value = '36.22'
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('D31').Value = '160.60'
$ws.Range('D32').Value = '5.51'
$ws.Range('E32').Value = '  +1.61%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '0.0770'
$ws.Range('E34').Value = '  +2.15%  '
$ws.Range('D36').Value = '17.20'
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').Value = '0.116'
$ws.Range('E38').Value = '  +1.63%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '1.83'
$ws.Range('E39').Value = '  +1.30%  '
$ws.Range('E40').Value = '  -3.65%  '
$ws.Range('D41').Value = '3.96'
$ws.Range('E41').Value = '  -2.71%  '
$ws.Range('E42').Value = '  +1.84%  '
$ws.Range('D43').Value = '1.991.37'
$ws.Range('E43').Value = '  +0.65%  '
$ws.Range('D44').Value = '19.12'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').Value = '0.0285'
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('E46').Value = '  +1.91%  '
$ws.Range('D47').Value = '9.19'
$ws.Range('E47').Value = '  +2.94%  '
$ws.Range('D48').Value = '2.712.39'
$ws.Range('E48').Value = '  -0.27%  '
$ws.Range('D49').Value = '97.61'
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').Value = '67.18'
$ws.Range('D51').Value = '73.51'
$ws.Range('E51').Value = '  +1.76%  '

# Restore default styling on the Price column (removes the temporary text format).
$ws.Range("D2:D51").Style = "Normal"

